# Generate Report for Handback
# Refresh the "Correspond Handoff Datetime" (E) and "Correspond Handback
# DateTime" (H) values for the ab3f5bcd... file row (row 2) on both the
# zh-cn and de-de language sheets, reflecting a new handback report run.
# The dcdaf9f5... file row (row 3) keeps its previously recorded timestamps.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("E2").Value = "2016-03-19 00:36:48"
$ws_zhcn.Range("H2").Value = "2016-03-19 00:37:08"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("E2").Value = "2016-03-19 00:36:50"
$ws_dede.Range("H2").Value = "2016-03-19 00:37:13"
